$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

# time_taken values for data rows (row 2..67)
$values = @(
    "2021-10-05 10:51:55.522527",
    "2021-10-05 10:51:55.522539",
    "2021-10-05 10:51:55.522542",
    "2021-10-05 10:51:55.522544",
    "2021-10-05 10:51:55.522547",
    "2021-10-05 10:51:55.522550",
    "2021-10-05 10:51:55.522552",
    "2021-10-05 10:51:55.522555",
    "2021-10-05 10:51:55.522558",
    "2021-10-05 10:51:55.522560",
    "2021-10-05 10:51:55.522563",
    "2021-10-05 10:51:55.522565",
    "2021-10-05 10:51:55.522567",
    "2021-10-05 10:51:55.522570",
    "2021-10-05 10:51:55.522572",
    "2021-10-05 10:51:55.522575",
    "2021-10-05 10:51:55.522578",
    "2021-10-05 10:51:55.522580",
    "2021-10-05 10:51:55.522583",
    "2021-10-05 10:51:55.522585",
    "2021-10-05 10:51:55.522587",
    "2021-10-05 10:51:55.522590",
    "2021-10-05 10:51:55.522592",
    "2021-10-05 10:51:55.522595",
    "2021-10-05 10:51:55.522598",
    "2021-10-05 10:51:55.522600",
    "2021-10-05 10:51:55.522603",
    "2021-10-05 10:51:55.522605",
    "2021-10-05 10:51:55.522607",
    "2021-10-05 10:51:55.522610",
    "2021-10-05 10:51:55.522612",
    "2021-10-05 10:51:55.522615",
    "2021-10-05 10:51:55.522617",
    "2021-10-05 10:51:55.522620",
    "2021-10-05 10:51:55.522622",
    "2021-10-05 10:51:55.522625",
    "2021-10-05 10:51:55.522627",
    "2021-10-05 10:51:55.522630",
    "2021-10-05 10:51:55.522633",
    "2021-10-05 10:51:55.522635",
    "2021-10-05 10:51:55.522638",
    "2021-10-05 10:51:55.522641",
    "2021-10-05 10:51:55.522643",
    "2021-10-05 10:51:55.522646",
    "2021-10-05 10:51:55.522648",
    "2021-10-05 10:51:55.522651",
    "2021-10-05 10:51:55.522653",
    "2021-10-05 10:51:55.522656",
    "2021-10-05 10:51:55.522658",
    "2021-10-05 10:51:55.522660",
    "2021-10-05 10:51:55.522663",
    "2021-10-05 10:51:55.522665",
    "2021-10-05 10:51:55.522668",
    "2021-10-05 10:51:55.522671",
    "2021-10-05 10:51:55.522673",
    "2021-10-05 10:51:55.522676",
    "2021-10-05 10:51:55.522678",
    "2021-10-05 10:51:55.522681",
    "2021-10-05 10:51:55.522683",
    "2021-10-05 10:51:55.522685",
    "2021-10-05 10:51:55.522688",
    "2021-10-05 10:51:55.522690",
    "2021-10-05 10:51:55.522693",
    "2021-10-05 10:51:55.522695",
    "2021-10-05 10:51:55.522699",
    "2021-10-05 10:51:55.522702"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
